# Auto-generated script applying FFXIV leve-profit market-price refresh
# as described by the commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2384.1538
$ws.Range("I2").Value = 715.7143
$ws.Range("J2").Value = 4330.6665
$ws.Range("K2").Value = 715.7143
$ws.Range("L2").Value = 4330.6665
$ws.Range("M2").Value = -602.7143
$ws.Range("N2").Value = -4556.6665

$ws.Range("H19").Value = 1574.5294
$ws.Range("I19").Value = 1594.7
$ws.Range("J19").Value = 1545.7142
$ws.Range("K19").Value = 1594.7
$ws.Range("L19").Value = 1545.7142
$ws.Range("M19").Value = -1419.7
$ws.Range("N19").Value = -1895.7142

$ws.Range("H28").Value = 183.2
$ws.Range("I28").Value = 183.2
$ws.Range("K28").Value = 183.2
$ws.Range("M28").Value = 301.8

$ws.Range("H43").Value = 750000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = $null

$ws.Range("H62").Value = 3694.8333
$ws.Range("I62").Value = 3500.4119
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 3500.4119
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -2876.4119
$ws.Range("N62").Value = -8248

$ws.Range("H65").Value = 3694.8333
$ws.Range("I65").Value = 3500.4119
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 17502.0595
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -14382.0595
$ws.Range("N65").Value = -41240

$ws.Range("H98").Value = 1691.3846
$ws.Range("J98").Value = 2679
$ws.Range("L98").Value = 2679
$ws.Range("N98").Value = -5675

$ws.Range("H116").Value = 5100.5
$ws.Range("I116").Value = 2216.5
$ws.Range("J116").Value = 7984.5
$ws.Range("K116").Value = 2216.5
$ws.Range("L116").Value = 7984.5
$ws.Range("M116").Value = 1225.5
$ws.Range("N116").Value = -14868.5

$ws.Range("H122").Value = 1691.3846
$ws.Range("J122").Value = 2679
$ws.Range("L122").Value = 8037
$ws.Range("N122").Value = -12937

$ws.Range("H125").Value = 3478.7693
$ws.Range("I125").Value = 2237.7144
$ws.Range("J125").Value = 4926.6665
$ws.Range("K125").Value = 20139.4296
$ws.Range("L125").Value = 44339.9985
$ws.Range("M125").Value = -17679.4296
$ws.Range("N125").Value = -49259.9985

$ws.Range("H132").Value = 3587.9092
$ws.Range("I132").Value = 3587.9092
$ws.Range("K132").Value = 10763.7276
$ws.Range("M132").Value = -8233.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2149.3333
$ws.Range("I2").Value = 2674
$ws.Range("J2").Value = 1100
$ws.Range("K2").Value = 2674
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = -2561
$ws.Range("N2").Value = -1326

$ws.Range("H4").Value = 606
$ws.Range("I4").Value = 749.5
$ws.Range("J4").Value = 32
$ws.Range("K4").Value = 749.5
$ws.Range("L4").Value = 32
$ws.Range("M4").Value = -633.5
$ws.Range("N4").Value = -264

$ws.Range("H32").Value = 2040.9846
$ws.Range("I32").Value = 1899.4286
$ws.Range("K32").Value = 1899.4286
$ws.Range("M32").Value = -1612.4286

$ws.Range("H116").Value = 2149.3333
$ws.Range("I116").Value = 2674
$ws.Range("J116").Value = 1100
$ws.Range("K116").Value = 2674
$ws.Range("L116").Value = 1100
$ws.Range("M116").Value = -380
$ws.Range("N116").Value = -5688

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2149.3333
$ws.Range("I3").Value = 2674
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 2674
$ws.Range("L3").Value = 1100
$ws.Range("M3").Value = -2560
$ws.Range("N3").Value = -1328

$ws.Range("H86").Value = 2342.7666
$ws.Range("I86").Value = 2554.3635
$ws.Range("K86").Value = 2554.3635
$ws.Range("M86").Value = -1431.3635

$ws.Range("H89").Value = 2342.7666
$ws.Range("I89").Value = 2554.3635
$ws.Range("K89").Value = 12771.8175
$ws.Range("M89").Value = -7155.817499999999

$ws.Range("H99").Value = 3463.8667
$ws.Range("I99").Value = 3477.5386
$ws.Range("K99").Value = 3477.5386
$ws.Range("M99").Value = -1979.5386

$ws.Range("H105").Value = 4006
$ws.Range("I105").Value = 3864.1428
$ws.Range("K105").Value = 3864.1428
$ws.Range("M105").Value = -2117.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = $null

$ws.Range("H31").Value = 10000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -10590

$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -10404

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null

$ws.Range("H134").Value = 2618.5
$ws.Range("I134").Value = 2618.5
$ws.Range("K134").Value = 7855.5
$ws.Range("M134").Value = -5320.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 5000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4884
$ws.Range("N3").Value = $null

$ws.Range("H126").Value = 9483.923000000001
$ws.Range("J126").Value = 6959.8
$ws.Range("L126").Value = 20879.4
$ws.Range("N126").Value = -25819.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 7667
$ws.Range("I2").Value = 6501
$ws.Range("J2").Value = 9999
$ws.Range("K2").Value = 6501
$ws.Range("L2").Value = 9999
$ws.Range("M2").Value = -6389
$ws.Range("N2").Value = -10223

$ws.Range("H100").Value = 773.44446
$ws.Range("I100").Value = 830.8333
$ws.Range("K100").Value = 1661.6666
$ws.Range("M100").Value = -1120.6666

$ws.Range("H122").Value = 4066.5
$ws.Range("I122").Value = 3879.8
$ws.Range("K122").Value = 11639.4
$ws.Range("M122").Value = -9189.400000000001
